$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Nombre) — values already exist in the shared-string table ---
$ws.Range("A1").Value = "Nombre"
$ws.Range("A2").Value = "Juan"
$ws.Range("A3").Value = "Luis"
$ws.Range("A4").Value = "Ana"

# --- Column B (new "Localizacion" column) ---
$ws.Range("B1").Value = "Localizacion"
$ws.Range("B2").Value = "105.70.56"
$ws.Range("B3").Value = "105.70.56"
$ws.Range("B4").Value = "109.8.9"

# --- Column C (email, unchanged, hyperlinked) ---
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("C2").Value = "juan@example.com"
$ws.Range("C3").Value = "luis@example.com"
$ws.Range("C4").Value = "ana@example.com"

# --- Column D (new "Identificador" column) ---
$ws.Range("D1").Value = "Identificador"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3

# --- Column E (new "Tipo" column) ---
$ws.Range("E1").Value = "Tipo"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1

# --- Drop the old Nacionalidad (F) / DNI (G) columns ---
$ws.Range("F1:F4").ClearContents()
$ws.Range("G1:G3").ClearContents()
$ws.Range("G4").Value = ""

# --- Selection / view state ---
$ws.Range("G7").Select()
